# Remove "mode" (Most common donation) row from each sheet of the
# analysis answer template. In each worksheet, row 6 holds the
# "Most common donation" label in column A; deleting the entire row
# shifts the remaining rows up and Excel prunes the now-unused shared
# string on save.

$wb = $excel.ActiveWorkbook

foreach ($name in @("Zip", "Weekday", "Month")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(6).Delete()
}

# Leave the selection/active sheet similar to the authored end state:
# Weekday and Zip end up with row 6 selected (the row that was just
# removed collapses onto the next row), and Month becomes the active
# sheet with the cursor parked below the data.
$wsZip = $wb.Worksheets.Item("Zip")
$wsZip.Range("A6:XFD6").Select()

$wsWeekday = $wb.Worksheets.Item("Weekday")
$wsWeekday.Range("A6:XFD6").Select()

$wsMonth = $wb.Worksheets.Item("Month")
$wsMonth.Activate()
$wsMonth.Range("A13").Select()
